$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46 - shifts existing rows 46:80 down to 47:81
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with the new daily price record
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 44596
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100103
$ws.Range("H46").Value = "Frutos de hueso (carozo)"
$ws.Range("I46").Value = 100103002
$ws.Range("J46").Value = "Ciruela"
$ws.Range("K46").Value = "Black Amber"
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 200
$ws.Range("N46").Value = 9000
$ws.Range("O46").Value = 9000
$ws.Range("P46").Value = 9000
$ws.Range("Q46").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R46").Value = "Provincia de Curicó"
$ws.Range("S46").Value = 500
$ws.Range("T46").Value = 18
